$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels from "god20XX" to "g20XX"
$ws.Range("C1").Value = "g2017"
$ws.Range("D1").Value = "g2018"
$ws.Range("E1").Value = "g2019"
$ws.Range("F1").Value = "g2020"
$ws.Range("G1").Value = "g2021"

# Update the active selection shown in the sheet view
$ws.Range("G4").Select()
